# Merge the "(x ", "= ", "1)" runs in the slide 9 title into a single
# run "(x = 1)", keeping the rest of the title text ("ПРФРА ") untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)          # "Заголовок 1" (title placeholder)

$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(7, 7)      # the "(x = 1)" portion of the title
$sub.Text = "(x = 1)"
